# Round the numeric data values in B2:E13 to the nearest integer.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B2:E13")

foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($null -ne $val) {
        $d = [double]$val
        if ($d -ge 0) {
            $rounded = [Math]::Floor($d + 0.5)
        } else {
            $rounded = [Math]::Ceiling($d - 0.5)
        }
        $cell.Value = $rounded
    }
}
